# Update column F ("dSF") values for the affected rows.
# Mapping of row -> new value (old values per the original sheet are
# overwritten with the repulled/recalculated data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    4  = 3
    6  = 2
    11 = 6
    13 = -2
    21 = -6
    22 = -3
    28 = 6
    40 = -2
    41 = -1
    43 = -2
    47 = 0
    50 = -1
    55 = 0
    57 = -2
    58 = -4
    59 = -11
    64 = -1
    73 = 6
    75 = 3
    77 = 3
    79 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
